$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Wealth Class in Allocation Row" -> "Current Income Row in IC Sheet", value 17 -> 4
$ws.Range("A2").Value = "Current Income Row in IC Sheet"
$ws.Range("B2").Value = 4

# Row 3: "Wealth Class in Cash Flow Row" -> "Current Expense Row in IC Sheet", value 20 -> 10
$ws.Range("A3").Value = "Current Expense Row in IC Sheet"
$ws.Range("B3").Value = 10

# Row 4: "Account Row" label unchanged, value 7 -> 6
$ws.Range("B4").Value = 6

# Row 5 "Wealth Row" / value 6 stays unchanged

# Old row 6 "Income Row" and row 7 "Expense Row" are removed, and old
# row 8 "Records Row" / row 9 "Records Banks Column" shift up to become
# the new rows 6 and 7.
$ws.Range("A6:B7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# New row 6: "Records Row", value 2 (unchanged content, now at row 6)
$ws.Range("A6").Value = "Records Row"
$ws.Range("B6").Value = 2

# New row 7: "Records Banks Column", value "G" (was "K")
$ws.Range("A7").Value = "Records Banks Column"
$ws.Range("B7").Value = "G"
